# "commit lanjutan untuk tutorial" - continuation pass over the slide
# transitions used in the Laravel tutorial deck.
#
# Slide 2 ("Kenapa Laravel" pillars diagram): Push-from-Left -> Pull-from-Top
# Slide 3: Push-from-Top -> Push-from-Left

$p = $ppt.ActivePresentation

# --- Slide 2: push(left) -> pull(up) -----------------------------------
$s2 = $p.Slides.Item(2)
$t2 = $s2.SlideShowTransition
$t2.EntryEffect     = 2050      # ppEffectPullUp
$t2.Speed           = 1         # ppTransitionSpeedSlow (spd="slow")
$t2.AdvanceOnClick  = -1        # msoTrue (advClick="1")

# --- Slide 3: push(up) -> push(left) ------------------------------------
$s3 = $p.Slides.Item(3)
$t3 = $s3.SlideShowTransition
$t3.EntryEffect     = 3853      # ppEffectPushLeft
$t3.Speed           = 1         # ppTransitionSpeedSlow (spd="slow")
$t3.AdvanceOnClick  = -1        # msoTrue (advClick="1")
